$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.536.32"
$ws.Range("E2").Value = "  -5.89%  "

# Row 3
$ws.Range("D3").Value = "2.979.15"
$ws.Range("E3").Value = "  -7.22%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.91"
$ws.Range("E5").Value = "  -6.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.76"
$ws.Range("E6").Value = "  -8.75%  "

# Row 7
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.560"
$ws.Range("E8").Value = "  -6.17%  "

# Row 9
$ws.Range("D9").Value = "2.985.57"
$ws.Range("E9").Value = "  -6.72%  "

# Row 10
$ws.Range("E10").Value = "  -7.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("E11").Value = "  -8.07%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("E12").Value = "  -7.03%  "

# Row 13
$ws.Range("D13").Value = "3.505.77"
$ws.Range("E13").Value = "  -6.87%  "

# Row 14
$ws.Range("E14").Value = "  -3.76%  "

# Row 15
$ws.Range("D15").Value = "61.711.27"
$ws.Range("E15").Value = "  -5.59%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.57"
$ws.Range("E16").Value = "  -8.40%  "

# Row 17
$ws.Range("D17").Value = "2.986.76"
$ws.Range("E17").Value = "  -6.65%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000146"
$ws.Range("E18").Value = "  -7.53%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "387.85"
$ws.Range("E19").Value = "  -6.28%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.09"
$ws.Range("E20").Value = "  -4.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.87"
$ws.Range("E21").Value = "  -8.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.59"
$ws.Range("E22").Value = "  -8.16%  "

# Row 23
$ws.Range("E23").Value = "  -0.14%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.74"
$ws.Range("E24").Value = "  -7.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.466"
$ws.Range("E25").Value = "  -4.92%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.185"
$ws.Range("E26").Value = "  -8.69%  "

# Row 27
$ws.Range("E27").Value = "  -0.47%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0930"
$ws.Range("E28").Value = "  -11.79%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.36"
$ws.Range("E29").Value = "  -6.48%  "

# Row 30
$ws.Range("E30").Value = "  -0.06%  "

# Row 31
$ws.Range("E31").Value = "  -7.76%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.24"
$ws.Range("E32").Value = "  -6.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.50"
$ws.Range("E33").Value = "  +1.26%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.98"
$ws.Range("E34").Value = "  -6.64%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.61"
$ws.Range("E35").Value = "  -8.27%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  -7.02%  "

# Row 37
$ws.Range("E37").Value = "  -7.09%  "

# Row 38
$ws.Range("E38").Value = "  -9.84%  "

# Row 39
$ws.Range("D39").Value = "2.429.62"
$ws.Range("E39").Value = "  -11.18%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.14"
$ws.Range("E40").Value = "  -5.19%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.87"
$ws.Range("E41").Value = "  -6.92%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.21"
$ws.Range("E42").Value = "  -8.32%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.657"
$ws.Range("E43").Value = "  -8.35%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0592"
$ws.Range("E44").Value = "  -6.74%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0244"
$ws.Range("E46").Value = "  -7.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.88"
$ws.Range("E47").Value = "  -12.62%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0951"
$ws.Range("E48").Value = "  -4.02%  "

# Row 49
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.48"
$ws.Range("E49").Value = "  +0.23%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.54"
$ws.Range("E50").Value = "  -9.56%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "263.28"
$ws.Range("E51").Value = "  -11.59%  "
